# 6023_W04_Threading_Notes.pptx
#
# Target change (slide 4, "Tools and tips"): the bullet
#   "Think about the DATA you are protecting, NOT the code you are running"
# gets a yellow highlight applied to its run.
#
# (The source diff also shows the Notes Master's auto date field cache
#  ticking from 2024-02-06 to 2024-02-12 - that's PowerPoint silently
#  re-caching the `datetimeFigureOut` field on open/save on a different
#  day, not a user edit, so there is nothing to author here.)

$p = $ppt.ActivePresentation

# Slide 4 = "Tools and tips"
$slide = $p.Slides.Item(4)

# The body placeholder holding the bullet list.
$shape = $slide.Shapes.Item(2)
$textRange = $shape.TextFrame.TextRange

$target = $textRange.Find("Think about the DATA you are protecting, NOT the code you are running")
if ($target -ne $null) {
    # RGB(255,255,0) = yellow, packed as 255 + 255*256 + 0*65536
    $target.Font.Highlight.RGB = 65535
}
